$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tmValues = @{
    2 = 220.0484447694182
    3 = 208.6084591231746
    4 = 213.8874033054168
    5 = 197.0219356245514
    6 = 170.5577719882461
    7 = 215.6270332331814
    8 = 216.0377001472469
    9 = 224.7239682442797
    10 = 256.0362566694457
    11 = 203.0371215863401
    12 = 280.4362278049025
    13 = 222.4008007752698
    14 = 165.078628675359
    15 = 165.4751946939818
    16 = 159.7187925647417
    17 = 169.0690506959194
    18 = 159.3250745690981
    19 = 153.3032293650398
    20 = 17.75289368952343
    21 = 179.2081797281337
    22 = 195.3409466607632
    23 = 183.9904554018607
    24 = 203.9507950247512
    25 = 300.1082255086671
    26 = 309.0042918694086
    27 = 257.7604894258689
    28 = 247.5140337670032
    29 = 392.589484603086
    30 = 246.2867492986718
    31 = 231.5488884934528
    32 = 260.6280371181025
    33 = 236.344307552025
    34 = 252.6990024029488
    35 = 254.971865362653
    36 = 207.563210835948
    37 = 196.1079124699839
    38 = 258.141103422786
    39 = 221.9777168520152
    40 = 202.0140417749864
    41 = 214.8559096480903
    42 = 376.5913657353018
    43 = 228.5526176208561
    44 = 230.8278103852844
    45 = 219.8785094082979
    46 = 182.2206999601753
    47 = 224.6826383798393
    48 = 250.3686568470001
    49 = 231.4414076049985
    50 = 296.3111446943747
    51 = 265.3037407571973
    52 = 306.2519379857882
    53 = 263.9503383484334
    54 = 235.2913552568154
    55 = 274.5816627951975
    56 = 288.9654023120442
    57 = 265.7244013424714
    58 = 246.2376294420228
    59 = 267.8943352054602
    60 = 235.1589854126764
    61 = 239.7088422903788
    62 = 259.0561610974572
    63 = 259.8588406962959
    64 = 243.0785237362668
    65 = 228.7488652083459
    66 = 254.744676274169
    67 = 231.5932228024334
    68 = 235.0626736655169
    69 = 275.5861042287879
    70 = 252.3717793954773
    71 = 244.2811314932216
    72 = 251.6362961653136
    73 = 237.3497831667671
    74 = 222.3517565342689
    75 = 224.6703814733575
    76 = 230.04601175482
    77 = 257.0125213383843
    78 = 266.8375267355315
    79 = 284.0219578714932
    80 = 282.0530583415607
    81 = 282.1638664087512
    82 = 289.1369392832293
    83 = 274.4483539737247
    84 = 256.3975647170884
    85 = 286.7781603204883
    86 = 299.2718719968973
    87 = 282.2945941102394
    88 = 210.1669720840187
    89 = 207.6897878591045
    90 = 213.6356742871829
    91 = 207.3289165846579
    92 = 206.0257248772071
    93 = 263.231846372648
    94 = 198.0157317683595
    95 = 242.7872843639836
    96 = 208.2873164481871
    97 = 219.0442014872023
    98 = 254.1227944175814
    99 = 244.1811711082867
    100 = 274.3150992006806
    101 = 239.3979208236884
    102 = 262.4713061076545
    103 = 226.4257692733164
    104 = 247.882956163346
    105 = 252.4903398476561
    106 = 250.0509020759595
    107 = 257.6665741797964
    108 = 270.6321236493634
    109 = 278.4430527142831
    110 = 318.7403632360842
    111 = 255.7802003854084
    112 = 276.2619745160586
    113 = 260.1117631476616
    114 = 433.0350720529276
    115 = 314.0605903920754
    116 = 308.4813173868908
    117 = 320.4141054859384
    118 = 325.5924095189797
    119 = 374.0247573730396
    120 = 375.2089721513517
    121 = 401.2614683952913
    122 = 319.2235760058658
    123 = 317.8105046517489
    124 = 207.6879945418398
    125 = 170.8362821838096
    126 = 199.8292598530225
    127 = 185.978606788332
    128 = 183.3162153773193
    129 = 224.1601527515164
    130 = 212.795426226645
    131 = 203.7735782838972
    132 = 203.0267808513176
    133 = 217.0133790067827
    134 = 343.4786997760781
    135 = 352.1522094027495
    136 = 300.9711189102819
    137 = 317.26489774442
    138 = 286.7979581781775
    139 = 428.2781753745896
    140 = 365.873033897836
    141 = 330.0427206108311
    142 = 352.9507984585313
    143 = 406.7101342679556
    144 = 232.5715969564982
    145 = 228.9865809680715
    146 = 258.4022846473553
}

foreach ($row in $tmValues.Keys) {
    $ws.Cells.Item([int]$row, 25).Value = $tmValues[$row]
}

Write-Output "Updated $($tmValues.Count) Y cells"